$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IMPORT_TNVED_6302")
$ws.Activate()

# The user selected A5:N6 (the two sample/demo rows that had been added
# for testing import) and deleted their contents, leaving the rows blank
# again but keeping whatever cell-level formatting was already present.
$rng = $ws.Range("A5:N6")
$rng.Select()
$rng.ClearContents()
